$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo / bad data values
$ws.Range("A2").Value = "2105670678_OLV"
$ws.Range("F2").Value = "[6,9]"
$ws.Range("E3").Value = "[700]"

# Row 2 no longer has the taller custom height (back to default row height)
$ws.Rows.Item(2).AutoFit()

# Update view: scroll so column D is the leftmost visible column, and select E3
$ws.Range("E3").Select()
$excel.ActiveWindow.ScrollColumn = 4
